$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P6").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("P8").Value = 285156
$ws.Range("P9").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("P11").Value = 0
$ws.Range("P12").Value = 186752
$ws.Range("P13").Value = 0
$ws.Range("P14").Value = 0
